$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.066364526748657
$ws.Range("B1").Value = 5.019322395324707
$ws.Range("C1").Value = 3.227057695388794
$ws.Range("D1").Value = 2.293357610702515
$ws.Range("E1").Value = 2.003292560577393
